$wb = $excel.ActiveWorkbook

# --- XPath sheet: remove the two rows that held the lone "call to customer"
# header (row 10) and the blank spacer above it (row 9). Deleting both shifts
# the remaining xpath rows up by two, matching the target layout and
# naturally drops the now-unused "call to customer" shared string.
$wsXPath = $wb.Worksheets.Item("XPath")
$wsXPath.Rows("9:10").Delete()

# --- Login sheet: zoom 110% -> 75%, keep existing selection (A2)
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Activate()
$wsLogin.Range("A2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75

# --- OrderInfo sheet: zoom 110% -> 75%, keep existing selection (A6)
$wsOrderInfo = $wb.Worksheets.Item("OrderInfo")
$wsOrderInfo.Activate()
$wsOrderInfo.Range("A6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75

# --- XPath sheet: zoom 110% -> 75%, new selection/top-left after the delete
$wsXPath.Activate()
$wsXPath.Range("B21").Select() | Out-Null
$excel.ActiveWindow.Zoom = 75
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
